$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6883953
$ws.Range("J17").Value = 6883953
$ws.Range("L17").Value = 20651859
$ws.Range("N17").Value = -20652195
$ws.Range("H32").Value = 6115.76
$ws.Range("I32").Value = 5367.5
$ws.Range("K32").Value = 5367.5
$ws.Range("M32").Value = -5041.5
$ws.Range("H100").Value = 1537.5
$ws.Range("I100").Value = 1633.6522
$ws.Range("K100").Value = 1633.6522
$ws.Range("M100").Value = -1092.6522
$ws.Range("H107").Value = 261.3125
$ws.Range("I107").Value = 295
$ws.Range("J107").Value = 160.25
$ws.Range("K107").Value = 295
$ws.Range("L107").Value = 160.25
$ws.Range("M107").Value = 1625
$ws.Range("N107").Value = -4000.25
$ws.Range("H113").Value = 37055430
$ws.Range("I113").Value = 50008656
$ws.Range("J113").Value = 46219.145
$ws.Range("K113").Value = 50008656
$ws.Range("L113").Value = 46219.145
$ws.Range("M113").Value = -50005402
$ws.Range("N113").Value = -52727.145
$ws.Range("H132").Value = 1318.8125
$ws.Range("I132").Value = 1335.5807
$ws.Range("K132").Value = 4006.7421
$ws.Range("M132").Value = -1476.7421
$ws.Range("H137").Value = 4225.543
$ws.Range("I137").Value = 4115.6772
$ws.Range("K137").Value = 12347.0316
$ws.Range("M137").Value = -9797.0316

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10549.412
$ws.Range("I32").Value = 10163.674
$ws.Range("K32").Value = 10163.674
$ws.Range("M32").Value = -9876.674000000001
$ws.Range("H61").Value = 9668.727999999999
$ws.Range("I61").Value = 10205.8
$ws.Range("J61").Value = 4298
$ws.Range("K61").Value = 10205.8
$ws.Range("L61").Value = 4298
$ws.Range("M61").Value = -9993.799999999999
$ws.Range("N61").Value = -4722
$ws.Range("H74").Value = 3342.7273
$ws.Range("I74").Value = 2575.5
$ws.Range("K74").Value = 2575.5
$ws.Range("M74").Value = -1701.5
$ws.Range("H77").Value = 3342.7273
$ws.Range("I77").Value = 2575.5
$ws.Range("K77").Value = 12877.5
$ws.Range("M77").Value = -8509.5
$ws.Range("H88").Value = 2948.6875
$ws.Range("I88").Value = 2752
$ws.Range("K88").Value = 2752
$ws.Range("M88").Value = -2346
$ws.Range("H91").Value = 2948.6875
$ws.Range("I91").Value = 2752
$ws.Range("K91").Value = 2752
$ws.Range("M91").Value = -1348
$ws.Range("H97").Value = 1084.8182
$ws.Range("I97").Value = 1132.1904
$ws.Range("J97").Value = 90
$ws.Range("K97").Value = 1132.1904
$ws.Range("L97").Value = 90
$ws.Range("M97").Value = -636.1904
$ws.Range("N97").Value = -1082
$ws.Range("H122").Value = 4399
$ws.Range("I122").Value = 3998.6667
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 11996.0001
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -9546.000100000001
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 2173.25
$ws.Range("I132").Value = 2173.25
$ws.Range("K132").Value = 6519.75
$ws.Range("M132").Value = -3989.75
$ws.Range("H136").Value = 9668.727999999999
$ws.Range("I136").Value = 10205.8
$ws.Range("J136").Value = 4298
$ws.Range("K136").Value = 30617.4
$ws.Range("L136").Value = 12894
$ws.Range("M136").Value = -28067.4
$ws.Range("N136").Value = -17994

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1622.96
$ws.Range("I94").Value = 1712.3478
$ws.Range("J94").Value = 595
$ws.Range("K94").Value = 1712.3478
$ws.Range("L94").Value = 595
$ws.Range("M94").Value = -1261.3478
$ws.Range("N94").Value = -1497
$ws.Range("H99").Value = 4223.727
$ws.Range("I99").Value = 3023.1667
$ws.Range("J99").Value = 5664.4
$ws.Range("K99").Value = 3023.1667
$ws.Range("L99").Value = 5664.4
$ws.Range("M99").Value = -1525.1667
$ws.Range("N99").Value = -8660.4
$ws.Range("H107").Value = 1200.8182
$ws.Range("I107").Value = 1232.7894
$ws.Range("K107").Value = 1232.7894
$ws.Range("M107").Value = 687.2106000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 583.75
$ws.Range("I22").Value = 384
$ws.Range("J22").Value = 916.6667
$ws.Range("K22").Value = 384
$ws.Range("L22").Value = 916.6667
$ws.Range("M22").Value = -34
$ws.Range("N22").Value = -1616.6667
$ws.Range("H99").Value = 5335.625
$ws.Range("J99").Value = 5787.3
$ws.Range("L99").Value = 5787.3
$ws.Range("N99").Value = -8783.299999999999
$ws.Range("H107").Value = 1390.5454
$ws.Range("I107").Value = 548.1053000000001
$ws.Range("J107").Value = 6726
$ws.Range("K107").Value = 548.1053000000001
$ws.Range("L107").Value = 6726
$ws.Range("M107").Value = 1371.8947
$ws.Range("N107").Value = -10566
$ws.Range("H126").Value = 5335.625
$ws.Range("J126").Value = 5787.3
$ws.Range("L126").Value = 17361.9
$ws.Range("N126").Value = -22301.9
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = $null
$ws.Range("H141").Value = 183517.6
$ws.Range("J141").Value = 207673.47
$ws.Range("L141").Value = 207673.47
$ws.Range("N141").Value = -218033.47

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4125
$ws.Range("J80").Value = 4125
$ws.Range("L80").Value = 12375
$ws.Range("N80").Value = -14247
$ws.Range("H83").Value = 4125
$ws.Range("J83").Value = 4125
$ws.Range("L83").Value = 37125
$ws.Range("N83").Value = -46485

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 714.0645
$ws.Range("I97").Value = 653.85
$ws.Range("J97").Value = 823.5454999999999
$ws.Range("K97").Value = 653.85
$ws.Range("L97").Value = 823.5454999999999
$ws.Range("M97").Value = -157.85
$ws.Range("N97").Value = -1815.5455
$ws.Range("H102").Value = 2686.2263
$ws.Range("I102").Value = 2201.238
$ws.Range("K102").Value = 2201.238
$ws.Range("M102").Value = -579.2379999999998
$ws.Range("H107").Value = 2608420.5
$ws.Range("I107").Value = 4563236
$ws.Range("K107").Value = 4563236
$ws.Range("M107").Value = -4561316
$ws.Range("H109").Value = 75000
$ws.Range("J109").Value = 75000
$ws.Range("L109").Value = 75000
$ws.Range("N109").Value = -77080
$ws.Range("H122").Value = 5992.143
$ws.Range("I122").Value = 6408.0835
$ws.Range("K122").Value = 19224.2505
$ws.Range("M122").Value = -16774.2505
$ws.Range("H132").Value = 3956.1667
$ws.Range("I132").Value = 3617.5186
$ws.Range("K132").Value = 10852.5558
$ws.Range("M132").Value = -8322.5558

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1999
$ws.Range("I2").Value = 1998
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1998
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1886
$ws.Range("N2").Value = -2224
$ws.Range("H61").Value = 54652.79
$ws.Range("I61").Value = 57439.055
$ws.Range("K61").Value = 57439.055
$ws.Range("M61").Value = -57237.055
$ws.Range("H68").Value = 76181.64
$ws.Range("I68").Value = 3421
$ws.Range("J68").Value = 148942.28
$ws.Range("K68").Value = 3421
$ws.Range("L68").Value = 148942.28
$ws.Range("M68").Value = -2672
$ws.Range("N68").Value = -150440.28
$ws.Range("H71").Value = 76181.64
$ws.Range("I71").Value = 3421
$ws.Range("J71").Value = 148942.28
$ws.Range("K71").Value = 17105
$ws.Range("L71").Value = 744711.4
$ws.Range("M71").Value = -13361
$ws.Range("N71").Value = -752199.4
$ws.Range("H113").Value = 54652.79
$ws.Range("I113").Value = 57439.055
$ws.Range("K113").Value = 57439.055
$ws.Range("M113").Value = -55269.055
$ws.Range("H122").Value = 5924.615
$ws.Range("I122").Value = 5638.1816
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 16914.5448
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -14464.5448
$ws.Range("N122").Value = -27400
$ws.Range("H132").Value = 267211.56
$ws.Range("I132").Value = 344648.53
$ws.Range("J132").Value = 9088.333000000001
$ws.Range("K132").Value = 1033945.59
$ws.Range("L132").Value = 27264.999
$ws.Range("M132").Value = -1031415.59
$ws.Range("N132").Value = -32324.999
$ws.Range("H136").Value = 64523892
$ws.Range("I136").Value = 37045290
$ws.Range("K136").Value = 111135870
$ws.Range("M136").Value = -111133320

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4285.7144
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 4166.6665
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 4166.6665
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -5414.6665
$ws.Range("H65").Value = 4285.7144
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 4166.6665
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 20833.3325
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -27073.3325
$ws.Range("H96").Value = 3720.6316
$ws.Range("I96").Value = 3903.125
$ws.Range("J96").Value = 3587.9092
$ws.Range("K96").Value = 3903.125
$ws.Range("L96").Value = 3587.9092
$ws.Range("M96").Value = -2530.125
$ws.Range("N96").Value = -6333.9092
$ws.Range("H107").Value = 1719.08
$ws.Range("I107").Value = 954.9375
$ws.Range("K107").Value = 2864.8125
$ws.Range("M107").Value = -944.8125
$ws.Range("H122").Value = 3905.375
$ws.Range("I122").Value = 2812.25
$ws.Range("K122").Value = 8436.75
$ws.Range("M122").Value = -5986.75
$ws.Range("H132").Value = 195717.4
$ws.Range("I132").Value = 241283.5
$ws.Range("K132").Value = 723850.5
$ws.Range("M132").Value = -721320.5
